# Update the sprint project tracker:
#  - Row 15 becomes the "Model Resnet50" task (dates + developer + Done status)
#  - Row 16 becomes the new "Model VGG19 dengan Imporvisasi FC Layer" task
#    (the old Resnet50/DenseNet169 evaluation rows are gone)
#  - Rows 17-20 shift up to the remaining Dev/Deploy tasks
#  - Rows 21-23 no longer have matching tasks, so their data is cleared

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring C:G of rows 15 and 16 to the same "filled data" formatting used by
# the rest of the task rows (thin-border date columns + green "Done" status)
# before writing their new values.
$ws.Range("C9:G9").Copy()
$ws.Range("C15:G15").PasteSpecial(-4122)
$ws.Range("C9:G9").Copy()
$ws.Range("C16:G16").PasteSpecial(-4122)

$ws.Range("B15").Value = "Model Resnet50"
$ws.Range("C15").Value = "18-12-2021"
$ws.Range("D15").Value = "19-12-2021"
$ws.Range("E15").Value = "19-12-2021"
$ws.Range("F15").Value = "Tiara Intana Sari"
$ws.Range("G15").Value = "Done"

$ws.Range("B16").Value = "Model VGG19 dengan Imporvisasi FC Layer"
$ws.Range("C16").Value = "16-12-2021"
$ws.Range("D16").Value = "18-12-2021"
$ws.Range("E16").Value = "18-12-2021"
$ws.Range("F16").Value = "Ulfah Nur Oktaviana"
$ws.Range("G16").Value = "Done"

$ws.Range("B17").Value = "Development Front-End Web"
$ws.Range("F17").Value = "Colaborasi"

$ws.Range("B18").Value = "Configurasi Flask"
$ws.Range("F18").Value = "Colaborasi"

$ws.Range("B19").Value = "Development Back-End Web"
$ws.Range("F19").Value = "Colaborasi"

$ws.Range("B20").Value = "Deployment "
$ws.Range("F20").Value = "Colaborasi"

$ws.Range("B21:F21").Clear()
$ws.Range("B22:F22").Clear()
$ws.Range("B23:F23").Clear()

$ws.Range("J11").Select()
